$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header rename: "NO." -> "count"
$ws.Range("B1").Value = "count"

# New row 13: an empty data cell that keeps the pre-update (general) cell
# formatting rather than the new percentage number format applied below.
$ws.Range("B12").Copy()
$ws.Range("B13").PasteSpecial(-4122)

# Updated data values (fractions) for rows 2-12
$values = @(
    0.1288022813688213,
    0.061874806641229241,
    0.065074510314309886,
    0.074430338882229691,
    0.067495739331454696,
    0.014218009478672985,
    0.09648518884721391,
    0.083189902467010898,
    0.027893030229071511,
    0.014122645969826661,
    0.023596553805630249
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Apply the new 4-decimal number format to the updated data range
$ws.Range("B2:B12").NumberFormat = "0.0000_ "
